$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: the ae38855b file (row 2) gets a freshly generated
# handoff/handback report timestamp pair.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-20 07:47:28"
$wsZhCn.Range("G2").Value = "2016-01-20 07:48:16"

# "de-de" sheet: the ae38855b file (row 2) gets a freshly generated
# handoff/handback report timestamp pair.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-20 07:47:40"
$wsDeDe.Range("G2").Value = "2016-01-20 07:48:34"
